$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" cell so it stays text even when the new value
# looks like a plain number (Excel would otherwise auto-convert it).
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "25.826.45"
$ws.Range("E2").Value = "  +11.61%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.735.60"
$ws.Range("E3").Value = "  +7.13%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "0.9988"
$ws.Range("E4").Value = "  +0.86%  "

# Row 5 - BNB
Set-TextCell "D5" "316.22"
$ws.Range("E5").Value = "  +4.88%  "

# Row 6 - USDC
Set-TextCell "D6" "0.9935"
$ws.Range("E6").Value = "  +1.13%  "

# Row 7 - XRP
Set-TextCell "D7" "0.3815"
$ws.Range("E7").Value = "  +4.13%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.3635"
$ws.Range("E8").Value = "  +6.48%  "

# Row 9 - OKB
Set-TextCell "D9" "50.64"
$ws.Range("E9").Value = "  +19.99%  "

# Row 10 - Polygon
Set-TextCell "D10" "1.224"
$ws.Range("E10").Value = "  +6.97%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.07691"
$ws.Range("E11").Value = "  +8.91%  "

# Row 12 - BinanceUSD
Set-TextCell "D12" "0.9952"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13 - Solana
Set-TextCell "D13" "21.74"
$ws.Range("E13").Value = "  +8.39%  "

# Row 14 - Polkadot
Set-TextCell "D14" "6.446"
$ws.Range("E14").Value = "  +9.37%  "

# Row 15 - Chainlink
Set-TextCell "D15" "7.061"
$ws.Range("E15").Value = "  +6.78%  "

# Row 16 - WrappedEther
Set-TextCell "D16" "1.745.04"
$ws.Range("E16").Value = "  +7.76%  "

# Row 17 - ShibaInu
Set-TextCell "D17" "0.00001154"
$ws.Range("E17").Value = "  +6.95%  "

# Row 18 - Dai (only E changes)
$ws.Range("E18").Value = "  +1.15%  "

# Row 19 - TRON
Set-TextCell "D19" "0.06819"
$ws.Range("E19").Value = "  +2.51%  "

# Row 20 - Litecoin
Set-TextCell "D20" "87.02"
$ws.Range("E20").Value = "  +11.18%  "

# Row 21 - Avalanche
Set-TextCell "D21" "17.55"
$ws.Range("E21").Value = "  +8.95%  "

# Row 22 - Uniswap
Set-TextCell "D22" "6.459"
$ws.Range("E22").Value = "  +7.76%  "

# Row 23 - Cosmos (only E changes)
$ws.Range("E23").Value = "  +9.31%  "

# Row 24 - WrappedBTC
Set-TextCell "D24" "25.764.99"
$ws.Range("E24").Value = "  +11.46%  "

# Row 25 - Toncoin
Set-TextCell "D25" "2.427"
$ws.Range("E25").Value = "  +2.22%  "

# Row 26 - LidoDAOToken
Set-TextCell "D26" "2.931"
$ws.Range("E26").Value = "  +12.99%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "20.63"
$ws.Range("E27").Value = "  +6.61%  "

# Row 28 - Monero
Set-TextCell "D28" "154.26"
$ws.Range("E28").Value = "  +2.90%  "

# Row 29 - BitcoinCash
Set-TextCell "D29" "134.15"
$ws.Range("E29").Value = "  +7.51%  "

# Row 30 - WrappedliquidstakedEther2.0
Set-TextCell "D30" "1.938.53"
$ws.Range("E30").Value = "  +7.95%  "

# Row 31 - ImmutableX
Set-TextCell "D31" "1.196"
$ws.Range("E31").Value = "  +22.99%  "

# Row 32 - Filecoin
Set-TextCell "D32" "7.016"
$ws.Range("E32").Value = "  +16.53%  "

# Row 33 - HuobiToken
Set-TextCell "D33" "4.373"
$ws.Range("E33").Value = "  +8.46%  "

# Row 34 - Aptos
Set-TextCell "D34" "14.26"
$ws.Range("E34").Value = "  +20.61%  "

# Row 35 - WEMIXTOKEN
Set-TextCell "D35" "1.797"
$ws.Range("E35").Value = "  +8.21%  "

# Row 36 - Stellar
Set-TextCell "D36" "0.08678"
$ws.Range("E36").Value = "  +5.51%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextCell "D37" "5.637"
$ws.Range("E37").Value = "  +8.85%  "

# Row 38 - Hedera
Set-TextCell "D38" "0.06722"
$ws.Range("E38").Value = "  +9.66%  "

# Row 39 - was VeChain, now FraxShare
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D39" "9.283"
$ws.Range("E39").Value = "  +6.56%  "

# Row 40 - was FraxShare, now VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D40" "0.02460"
$ws.Range("E40").Value = "  +11.00%  "

# Row 41 - was TrustWalletToken, now Algorand
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D41" "0.2214"
$ws.Range("E41").Value = "  +9.64%  "

# Row 42 - was Algorand, now TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D42" "1.298"
$ws.Range("E42").Value = "  +3.70%  "

# Row 43 - TheSandbox (only D/E change)
Set-TextCell "D43" "0.6545"
$ws.Range("E43").Value = "  +11.14%  "

# Row 44 - was Frax, now EnergySwap
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D44" "13.93"
$ws.Range("E44").Value = "  +6.67%  "

# Row 45 - was EnergySwap, now Frax
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D45" "0.9920"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46 - Decentraland
Set-TextCell "D46" "0.6336"
$ws.Range("E46").Value = "  +10.40%  "

# Row 47 - PancakeSwap
Set-TextCell "D47" "3.900"
$ws.Range("E47").Value = "  +3.22%  "

# Row 48 - NEARProtocol
Set-TextCell "D48" "2.174"
$ws.Range("E48").Value = "  +10.13%  "

# Row 49 - Quant
Set-TextCell "D49" "131.82"
$ws.Range("E49").Value = "  +4.63%  "

# Row 50 - Cronos
Set-TextCell "D50" "0.07450"
$ws.Range("E50").Value = "  +7.42%  "

# Row 51 - Aave
Set-TextCell "D51" "79.25"
$ws.Range("E51").Value = "  +7.30%  "
